# Generate Report for Handback
# - Status text updated everywhere it appears (Overview + per-locale sheets)
#   from "Handed back: in sync with en-US" to "Handed back: not in sync with en-US".
# - de-de sheet: the handback round for the second file (e49630ad-...) produced a
#   fresh handback datetime, so "de-de"!K2 moves forward.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep $oldStatus on the left of -eq. Some cells (e.g. the
        # "True"/"False" text cells) come back from Value2 as native
        # booleans; with the boolean on the left, -eq would coerce this
        # string to a boolean and (mis)match every non-empty string.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# Refresh the Correspond Handback DateTime for the de-de handback of
# e49630ad-40ac-44e9-a236-cf67916150cf (row 2 in the de-de table).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-21 01:12:50"

# The status column/columns got visibly wider now that the longer
# "not in sync" text no longer fits the previous auto-fit width.
$newWidth = 32.666666666666664
$wb.Worksheets.Item("Overview").Range("E1").ColumnWidth = $newWidth
$wb.Worksheets.Item("Overview").Range("F1").ColumnWidth = $newWidth
$wb.Worksheets.Item("zh-cn").Range("C1").ColumnWidth = $newWidth
$wb.Worksheets.Item("de-de").Range("C1").ColumnWidth = $newWidth
